$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new column D values ("U") for data rows 2-6, and clear out the old
# F/G columns that used to hold the sire/dam numeric codes.
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 4).Value = "U"
    $ws.Cells.Item($r, 6).ClearContents()
    $ws.Cells.Item($r, 7).ClearContents()
}

# Update the active selection to match the author's saved cursor position.
$ws.Range("D7").Select()
